$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106..154 down to 107..155.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new weekly data point.
$ws.Cells.Item(106, 1).Value = 9
$ws.Cells.Item(106, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(106, 3).Value = "Metropolitana"
$ws.Cells.Item(106, 4).Value = 45006
$ws.Cells.Item(106, 5).Value = 13
$ws.Cells.Item(106, 6).Value = 100112022
$ws.Cells.Item(106, 7).Value = "Arveja Verde"
$ws.Cells.Item(106, 8).Value = "Perfection"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 52
$ws.Cells.Item(106, 11).Value = 27000
$ws.Cells.Item(106, 12).Value = 30000
$ws.Cells.Item(106, 13).Value = 28500
$ws.Cells.Item(106, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(106, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(106, 16).Value = 1140
$ws.Cells.Item(106, 17).Value = 25
$ws.Cells.Item(106, 18).Value = "Hortaliza"
